# Updates the cryptos list (price/volume refresh + two pairs of rows whose
# coin/link/price/volume got re-ordered) per the GitHub Actions data refresh.
#
# Note: several "Price" values are plain decimal-looking numbers (e.g.
# "262.78", "108.50", "0.06920") that must stay as literal text, exactly as
# authored upstream, instead of being auto-coerced to numbers (which would
# silently drop meaningful trailing zeros, e.g. "108.50" -> 108.5). A leading
# apostrophe forces Excel to keep them as text, same as typing them in
# manually.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---------------------------------------------------------
$ws.Range("D2").Value = "26.613.88"
$ws.Range("E2").Value = "  -0.12%  "

# --- Row 3: Ethereum ---------------------------------------------------------
$ws.Range("D3").Value = "1.849.55"
$ws.Range("E3").Value = "  -0.29%  "

# --- Row 5: BNB ---------------------------------------------------------
$ws.Range("D5").Value = "'262.78"
$ws.Range("E5").Value = "  -0.81%  "

# --- Row 6: USDC (volume only) ---------------------------------------------------------
$ws.Range("E6").Value = "  +0.12%  "

# --- Row 7: XRP (volume only) ---------------------------------------------------------
$ws.Range("E7").Value = "  +2.14%  "

# --- Row 8: Cardano ---------------------------------------------------------
$ws.Range("D8").Value = "'0.3146"
$ws.Range("E8").Value = "  -4.30%  "

# --- Row 9: Dogecoin ---------------------------------------------------------
$ws.Range("D9").Value = "'0.06920"
$ws.Range("E9").Value = "  +1.73%  "

# --- Row 10: Solana ---------------------------------------------------------
$ws.Range("D10").Value = "'18.82"
$ws.Range("E10").Value = "  -0.25%  "

# --- Row 11: Polygon ---------------------------------------------------------
$ws.Range("D11").Value = "'0.7694"
$ws.Range("E11").Value = "  -1.03%  "

# --- Row 12: TRON ---------------------------------------------------------
$ws.Range("D12").Value = "'0.07831"
$ws.Range("E12").Value = "  +1.36%  "

# --- Row 13: WrappedEther ---------------------------------------------------------
$ws.Range("D13").Value = "1.848.95"
$ws.Range("E13").Value = "  -0.27%  "

# --- Row 14: Litecoin ---------------------------------------------------------
$ws.Range("D14").Value = "'89.63"
$ws.Range("E14").Value = "  +1.08%  "

# --- Row 15: Polkadot ---------------------------------------------------------
$ws.Range("D15").Value = "'5.051"
$ws.Range("E15").Value = "  +0.23%  "

# --- Rows 16/17: Avalanche & BinanceUSD swap order, with refreshed values ---
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'14.10"
$ws.Range("E17").Value = "  +0.55%  "

# --- Row 18: ShibaInu ---------------------------------------------------------
$ws.Range("D18").Value = "'0.000007962"
$ws.Range("E18").Value = "  -0.37%  "

# --- Row 20: WrappedBTC ---------------------------------------------------------
$ws.Range("D20").Value = "26.634.24"
$ws.Range("E20").Value = "  -0.09%  "

# --- Row 21: WrappedliquidstakedEther2.0 ---------------------------------------------------------
$ws.Range("D21").Value = "2.085.71"
$ws.Range("E21").Value = "  -0.05%  "

# --- Row 22: Uniswap ---------------------------------------------------------
$ws.Range("D22").Value = "'4.646"
$ws.Range("E22").Value = "  +0.11%  "

# --- Row 23: Chainlink ---------------------------------------------------------
$ws.Range("D23").Value = "'6.019"
$ws.Range("E23").Value = "  +0.18%  "

# --- Row 24: Cosmos ---------------------------------------------------------
$ws.Range("D24").Value = "'9.344"
$ws.Range("E24").Value = "  -2.19%  "

# --- Row 25: LidoDAOToken ---------------------------------------------------------
$ws.Range("D25").Value = "'2.215"
$ws.Range("E25").Value = "  +0.69%  "

# --- Row 26: Monero ---------------------------------------------------------
$ws.Range("D26").Value = "'141.42"
$ws.Range("E26").Value = "  -2.05%  "

# --- Row 27: Toncoin ---------------------------------------------------------
$ws.Range("D27").Value = "'1.691"
$ws.Range("E27").Value = "  +0.66%  "

# --- Row 28: EthereumClassic ---------------------------------------------------------
$ws.Range("D28").Value = "'17.03"
$ws.Range("E28").Value = "  +0.22%  "

# --- Row 29: BitcoinCash ---------------------------------------------------------
$ws.Range("D29").Value = "'111.51"
$ws.Range("E29").Value = "  -0.82%  "

# --- Row 30: InternetComputer(DFINITY) ---------------------------------------------------------
$ws.Range("D30").Value = "'4.303"
$ws.Range("E30").Value = "  +2.59%  "

# --- Row 31: Stellar ---------------------------------------------------------
$ws.Range("D31").Value = "'0.08778"
$ws.Range("E31").Value = "  +0.19%  "

# --- Row 32: Filecoin ---------------------------------------------------------
$ws.Range("D32").Value = "'4.112"
$ws.Range("E32").Value = "  -1.11%  "

# --- Row 33: Hedera (volume only) ---------------------------------------------------------
$ws.Range("E33").Value = "  +0.41%  "

# --- Row 34: ImmutableX ---------------------------------------------------------
$ws.Range("D34").Value = "'0.7383"
$ws.Range("E34").Value = "  +3.22%  "

# --- Rows 35/36: ARBITRUM & HuobiToken swap order, with refreshed values ---
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.913"
$ws.Range("E35").Value = "  +1.56%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.139"
$ws.Range("E36").Value = "  -0.01%  "

# --- Row 38: RenderToken ---------------------------------------------------------
$ws.Range("D38").Value = "'2.333"
$ws.Range("E38").Value = "  +6.13%  "

# --- Row 39: VeChain ---------------------------------------------------------
$ws.Range("D39").Value = "'0.01733"
$ws.Range("E39").Value = "  -2.95%  "

# --- Row 40: TheSandbox ---------------------------------------------------------
$ws.Range("D40").Value = "'0.4829"
$ws.Range("E40").Value = "  -1.26%  "

# --- Row 41: TrustWalletToken ---------------------------------------------------------
$ws.Range("D41").Value = "'0.9054"
$ws.Range("E41").Value = "  +0.43%  "

# --- Row 42: Quant ---------------------------------------------------------
$ws.Range("D42").Value = "'108.50"
$ws.Range("E42").Value = "  -4.00%  "

# --- Row 43: FraxShare ---------------------------------------------------------
$ws.Range("D43").Value = "'5.907"
$ws.Range("E43").Value = "  -2.86%  "

# --- Row 44: PaxDollar (volume only) ---------------------------------------------------------
$ws.Range("E44").Value = "  +0.12%  "

# --- Row 45: Aptos ---------------------------------------------------------
$ws.Range("D45").Value = "'7.688"
$ws.Range("E45").Value = "  -0.62%  "

# --- Row 46: Decentraland ---------------------------------------------------------
$ws.Range("D46").Value = "'0.4191"
$ws.Range("E46").Value = "  -0.37%  "

# --- Rows 47/48: Algorand & EnergySwap swap order, with refreshed values ---
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.097"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1248"
$ws.Range("E48").Value = "  +0.52%  "

# --- Row 49: Elrond ---------------------------------------------------------
$ws.Range("D49").Value = "'35.06"
$ws.Range("E49").Value = "  -0.06%  "

# --- Row 50: Cronos ---------------------------------------------------------
$ws.Range("D50").Value = "'0.05814"
$ws.Range("E50").Value = "  -1.84%  "

# --- Row 51: EOS ---------------------------------------------------------
$ws.Range("D51").Value = "'0.8960"
$ws.Range("E51").Value = "  +1.27%  "
